# Se cambió la estructura del campo número de planilla en el módulo de radicación
#
# Sheet "13_03_2019" tracks requirement-tickets with a "TIEMPO TOTAL" (total
# time) column. Row 2 is the ticket "Cambiar el capo de número de planilla de
# acuerdo a la estructura dada." ("Change the payroll-number field structure").
# Its total time is being corrected, and row 10 (the grand-total row, a plain
# manually maintained sum, not a worksheet formula) is updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("13_03_2019")

# TIEMPO TOTAL for the "número de planilla" requirement: was ~1h10m, now 4h.
$ws.Range("F2").Value = 0.16666666666666666

# Grand total (TOTAL HORAS row) = SUM(D10,E10) + F2; keep it consistent with
# the corrected F2 above.
$ws.Range("F10").Value = 1.2708333333333333

# Leave the cursor where the user clicked after making the edit.
$ws.Range("A4").Select()
